$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Workbook-level metadata
# ---------------------------------------------------------------------
$ws.Name = "ts"

# ---------------------------------------------------------------------
# Insert a new column at H, shifting the existing "address" data (and
# everything to its right) one column over to I, J, K, ...
# ---------------------------------------------------------------------
$ws.Columns("H").Insert()

# New header cell for the inserted column.
$ws.Range("H1").Value2 = "office type "

# Fill the "Headquarter" office-type value for every row that already had
# contact/address data in the old H column (i.e. every row except the
# pure-continuation rows 4, 5, 6, and the trailing row 25).
$hqRows = @(2, 3, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)
foreach ($r in $hqRows) {
    $ws.Range("H$r").Value2 = "Headquarter"
}

# ---------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------
$ws.Range("H4").Select()

Write-Output "done"
